$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "70.028.89"
$ws.Range("E2").Value = "  -1.37%  "

$ws.Range("D3").Value = "3.747.25"
$ws.Range("E3").Value = "  +2.35%  "

$ws.Range("E4").Value = "  +0.01%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "622.62"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.42%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "180.05"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.90%  "

$ws.Range("D7").Value = "3.746.35"
$ws.Range("E7").Value = "  +2.46%  "

$ws.Range("E9").Value = "  -1.46%  "

$ws.Range("E10").Value = "  +3.09%  "

$ws.Range("E11").Value = "  -5.32%  "

$ws.Range("E12").Value = "  -3.04%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "41.05"
$ws.Range("D13").Style = "Normal"

$ws.Range("E14").Value = "  +2.17%  "

$ws.Range("D15").Value = "4.361.29"
$ws.Range("E15").Value = "  +2.04%  "

$ws.Range("D16").Value = "3.747.61"
$ws.Range("E16").Value = "  +1.71%  "

$ws.Range("D17").Value = "70.045.27"
$ws.Range("E17").Value = "  -1.34%  "

$ws.Range("E18").Value = "  -1.29%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.60"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.83%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "16.84"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.41%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "506.32"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.64%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.39"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.76%  "

$ws.Range("E23").Value = "  -2.21%  "

$ws.Range("E24").Value = "  +0.39%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "86.84"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.91%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "13.12"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.73%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.21"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.08%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0000138"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +24.27%  "

$ws.Range("E29").Value = "  -0.01%  "

$ws.Range("E30").Value = "  -2.64%  "

$ws.Range("E31").Value = "  +0.35%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.90"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.38%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "31.28"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.95%  "

$ws.Range("E34").Value = "  -0.31%  "

$ws.Range("E35").Value = "  -0.05%  "

$ws.Range("E36").Value = "  +3.96%  "

$ws.Range("E37").Value = "  +1.49%  "

$ws.Range("E38").Value = "  -3.94%  "

$ws.Range("E39").Value = "  +1.36%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.12"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -4.07%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "50.41"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.58%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "45.10"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.77%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "424.73"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.64%  "

$ws.Range("E44").Value = "  -1.16%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.85"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.24%  "

$ws.Range("D46").Value = "2.997.19"
$ws.Range("E46").Value = "  -3.64%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0365"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.66%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "27.34"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.15%  "

$ws.Range("E49").Value = "  -0.03%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "138.00"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.21%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.51"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.39%  "
